# Weekly price-sheet update: insert a new daily record as row 505 on the
# single data sheet ("Hortaliza, Feria Lagunitas de Puerto Montt -
# Zanahoria"). Inserting the row pushes the existing rows 505:540 down to
# 506:541 (dimension grows from A1:R540 to A1:R541), matching the rest of
# the dataset which is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 505, shifting rows 505-540 down
# to 506-541.
$ws.Rows.Item(505).Insert()

# Populate the newly inserted row 505 with the new observation.
$ws.Range("A505").Value = 4
$ws.Range("B505").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C505").Value = "Los Lagos"
$ws.Range("D505").Value = 45021
$ws.Range("E505").Value = 10
$ws.Range("F505").Value = 100114013
$ws.Range("G505").Value = "Zanahoria"
$ws.Range("H505").Value = "Sin especificar"
$ws.Range("I505").Value = "Primera"
$ws.Range("J505").Value = 150
$ws.Range("K505").Value = 8000
$ws.Range("L505").Value = 8000
$ws.Range("M505").Value = 8000
$ws.Range("N505").Value = "$/saco 20 kilos"
$ws.Range("O505").Value = "Región de La Araucanía"
$ws.Range("P505").Value = 400
$ws.Range("Q505").Value = 20
$ws.Range("R505").Value = "Hortaliza"
